# update testcases Search Home
#
# 1) Add five new worksheets at the end of the workbook with their data,
#    bold headers, page setup and leftover cell selections.
# 2) Restore focus to the "Search" sheet and move its remembered selection
#    from D7 to A2.
# 3) Best-effort nudge of the tab scroll position (firstSheet).

$wb = $excel.ActiveWorkbook

function Add-SheetAfter($afterSheet, $name) {
    $s = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
    $s.Name = $name
    return $s
}

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

# ---------------------------------------------------------------------
# Menu_Pilar_Tampil
# ---------------------------------------------------------------------
$menuPilarTampil = Add-SheetAfter $lastSheet "Menu_Pilar_Tampil"
$menuPilarTampil.Range("A1").Value = "Menu Id"
$menuPilarTampil.Range("A1").Font.Bold = $true
$menuPilarTampil.Range("A2").Value = "action-home"
$menuPilarTampil.Range("A3").Value = "action-live-event"
$menuPilarTampil.Range("A4").Value = "action-library"
$menuPilarTampil.Columns("A:A").ColumnWidth = 14.7265625
$menuPilarTampil.PageSetup.PaperSize = 9
$menuPilarTampil.PageSetup.Orientation = 1
$menuPilarTampil.Range("K14").Select()

# ---------------------------------------------------------------------
# Menu_Pilar_Not
# ---------------------------------------------------------------------
$menuPilarNot = Add-SheetAfter $menuPilarTampil "Menu_Pilar_Not"
$menuPilarNot.Range("A1").Value = "Menu"
$menuPilarNot.Range("A1").Font.Bold = $true
$menuPilarNot.Range("A2").Value = "action-live-tv"
$menuPilarNot.Range("A3").Value = "action-account"
$menuPilarNot.Columns("A:A").ColumnWidth = 13.1796875
$menuPilarNot.Range("G4").Select()

# ---------------------------------------------------------------------
# Pilar
# ---------------------------------------------------------------------
$pilar = Add-SheetAfter $menuPilarNot "Pilar"
$pilar.Range("A1").Value = "Pilar"
$pilar.Range("A1").Font.Bold = $true
$pilar.Range("B1").Value = "Content"
$pilar.Range("B1").Font.Bold = $true
$pilar.Range("A2").Value = "Videos"
$pilar.Range("B2").Value = "add-tab-button"
$pilar.Range("A3").Value = "News"
$pilar.Range("B3").Value = "Tab Search"
$pilar.Range("A4").Value = "Radio +"
$pilar.Range("B4").Value = "Program"
$pilar.Range("A5").Value = "Home of Talent"
$pilar.Range("B5").Value = "Episode"
$pilar.Range("A6").Value = "Games"
$pilar.Range("B6").Value = "Extra"
$pilar.Columns("A:A").ColumnWidth = 13.7265625
$pilar.Columns("B:B").ColumnWidth = 21.6328125
$pilar.PageSetup.PaperSize = 9
$pilar.PageSetup.Orientation = 1
$pilar.Range("H11").Select()

# ---------------------------------------------------------------------
# News_Content
# ---------------------------------------------------------------------
$newsContent = Add-SheetAfter $pilar "News_Content"
$newsContent.Range("A1").Value = "Content"
$newsContent.Range("A1").Font.Bold = $true
$newsContent.Range("A2").Value = "article-thumbnail"
$newsContent.Range("A3").Value = "add-tab-button"
$newsContent.Columns("A:A").ColumnWidth = 15.1796875
$newsContent.PageSetup.PaperSize = 9
$newsContent.PageSetup.Orientation = 1
$newsContent.Range("D7").Select()

# ---------------------------------------------------------------------
# Tab Search
# ---------------------------------------------------------------------
$tabSearch = Add-SheetAfter $newsContent "Tab Search"
$tabSearch.Range("A1").Value = "Index"
$tabSearch.Range("A1").Font.Bold = $true
$tabSearch.Range("B1").Value = "Direct"
$tabSearch.Range("B1").Font.Bold = $true
$tabSearch.Range("A2").Value = "0"
$tabSearch.Range("B2").Value = "/"
$tabSearch.Range("A3").Value = "1"
$tabSearch.Range("B3").Value = "/trending"
$tabSearch.Range("A4").Value = "2"
$tabSearch.Range("B4").Value = "/radio"
$tabSearch.Range("A5").Value = "3"
$tabSearch.Range("B5").Value = "Clip"
$tabSearch.Range("A6").Value = "4"
$tabSearch.Range("B6").Value = "Photo"
$tabSearch.Columns("A:B").ColumnWidth = 9.90625
$tabSearch.PageSetup.PaperSize = 9
$tabSearch.PageSetup.Orientation = 1
$tabSearch.Range("F13").Select()

# ---------------------------------------------------------------------
# Restore the active sheet/selection on "Search" (D7 -> A2) and nudge the
# tab scroll position back towards where it was before (firstSheet 30 -> 27).
# ---------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollWorkbookTabs(1, 27) | Out-Null

$search = $wb.Worksheets.Item("Search")
$search.Activate()
$search.Range("A2").Select()
